$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column N mirrors column M's styling (copy M1:M20 -> N1:N20), then set the
# new/updated values for the T12 (27/3/2020) snapshot.
$ws.Range("M1:M20").Copy()
$ws.Range("N1:N20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("N1").Value = "T12: 27/3/2020"

$ws.Range("N2").Value = 1
$ws.Range("N3").Value = 2
$ws.Range("N4").Value = 12
$ws.Range("N5").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("N7").Value = 43
$ws.Range("N8").Value = 0
$ws.Range("N9").Value = 32
$ws.Range("N10").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("N14").Value = 3
$ws.Range("N15").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("N19").Value = 2

$ws.Range("N20").Formula = "=SUM(N2:N19)"

# Column width for new column N (target stored width 16.39 chars; engine quantizes
# ColumnWidth to 1/6-character pixel steps, so 15.4167 is the closest achievable input)
$ws.Columns.Item(14).ColumnWidth = 15.416666666666666

# Right-align the B20/C20 total cells (new style introduced)
$ws.Range("B20:C20").HorizontalAlignment = -4152

# Update selection to D13
$ws.Range("D13").Select()
